$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Cost" and "Currency" columns (E, F) with per-row amounts, all in HUF.
# Written before the "Ticker" header so the shared-string table grows in the
# same order as the target workbook (Cost, Currency, HUF, then Ticker).
$ws.Range("E1").Value = "Cost"
$ws.Range("F1").Value = "Currency"

$costs = @(250, 300, 500, 750)
for ($i = 0; $i -lt $costs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $costs[$i]
    $ws.Cells.Item($row, 6).Value = "HUF"
}

# Rename the "Instrument" header to "Ticker".
$ws.Range("B1").Value = "Ticker"

# Match the saved selection state from the authored workbook.
$ws.Range("B2").Select()
